$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin rows (B:Coin, C:Link, D:Price, E:Volume(1h)) to latest scraped values
$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.117.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.076.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.95%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.61%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +16.85%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.351"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.52%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.074.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.81%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.716"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.93%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.195"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.72%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.81%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.769.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.46%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.630.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.88%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.059.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.60%  "

$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.58%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("B21").Value = "PEPE"
$ws.Range("C21").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000204"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.53%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "445.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.90%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("B24").Value = "Polkadot"
$ws.Range("C24").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.57%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.92%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.21%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.53%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.222.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.79%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.53%  "

$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.159"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.46%  "

$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.996"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.20%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +18.91%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.201"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +40.40%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.150"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.33%  "

$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "496.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.12%  "

$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.71%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -11.22%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.425"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.45%  "

$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0850"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.10%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("B46").Value = "MantraDAO"
$ws.Range("C46").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.44%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "147.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.81%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.686"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.88%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.36%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.92%  "

$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.13%  "

